$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between the two detail rows (2207 <-> 2206)
$ws.Range("E16").Value = "2206"
$ws.Range("E17").Value = "2207"

# Update the "Valor Mora" amounts for both rows
$ws.Range("G16").Value = 908526
$ws.Range("G17").Value = 908526
